$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 74; this shifts the existing rows 74-152
# down to 75-153 and extends the sheet dimension automatically.
$ws.Rows("74:74").Insert()

# Populate the newly inserted row 74 with the new weekly price record.
$ws.Range("A74").Value = 10
$ws.Range("B74").Value = "Vega Modelo de Temuco"
$ws.Range("C74").Value = "La Araucanía"
$ws.Range("D74").Value = 45033
$ws.Range("D74").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E74").Value = 9
$ws.Range("F74").Value = 100112035
$ws.Range("G74").Value = "Bruselas (repollito)"
$ws.Range("H74").Value = "Sin especificar"
$ws.Range("I74").Value = "Primera"
$ws.Range("J74").Value = 85
$ws.Range("K74").Value = 30000
$ws.Range("L74").Value = 30000
$ws.Range("M74").Value = 30000
$ws.Range("N74").Value = '$/malla 15 kilos'
$ws.Range("O74").Value = "Región Metropolitana"
$ws.Range("P74").Value = 2000
$ws.Range("Q74").Value = 15
$ws.Range("R74").Value = "Hortaliza"
